$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected cell on the sheet (was A29, now D20)
$ws.Range("D20").Select()

# Original error values per row (rows 2..19) for columns B (Train Error),
# C (Valid Error) and D (Test Error). Each cell is converted from the raw
# error value into a "1 - error" formula (i.e. accuracy instead of error).
$values = @{
    "B2"  = "0.00114351057747"
    "C2"  = "0.0108695652174"
    "D2"  = "0.0502702702703"
    "B3"  = "0.00664451827243"
    "C3"  = "0.027027027027"
    "D3"  = "0.0489887640449"
    "B4"  = "0.00069060773480667"
    "C4"  = "0.0526315789473685"
    "D4"  = "0.0606653620352251"
    "B5"  = "0.00190566936636"
    "C5"  = "0.0272727272727"
    "D5"  = "0.0691056910569"
    "B6"  = "0"
    "C6"  = "0"
    "D6"  = "0.0460829493088"
    "B7"  = "0.000473709142586487"
    "C7"  = "0.036036036036036"
    "D7"  = "0.0613718411552346"
    "B8"  = "0"
    "C8"  = "0.0833333333333"
    "D8"  = "0.0857142857143"
    "B9"  = "0"
    "C9"  = "0.0180180180180181"
    "D9"  = "0.0470588235294118"
    "B10" = "0.0159726183685"
    "C10" = "0.0108695652174"
    "D10" = "0.0573903627504"
    "B11" = "0.00996677740864"
    "C11" = "0.0540540540541"
    "D11" = "0.0572587917042"
    "B12" = "0"
    "C12" = "0.0833333333333"
    "D12" = "0.065306122449"
    "B13" = "0.00709219858156"
    "C13" = "0"
    "D13" = "0.0753424657534"
    "B14" = "0.00429184549356"
    "C14" = "0"
    "D14" = "0.0491803278689"
    "B15" = "0"
    "C15" = "0"
    "D15" = "0.0456621004566"
    "B16" = "0.093131548312"
    "C16" = "0.121546961326"
    "D16" = "0.105918141593"
    "B17" = "0.00189753320683"
    "C17" = "0.0630630630631"
    "D17" = "0.0803974706414"
    "B18" = "0"
    "C18" = "0"
    "D18" = "0.0462962962963"
    "B19" = "0"
    "C19" = "0.181818181818"
    "D19" = "0.0601851851852"
}

foreach ($addr in $values.Keys) {
    $orig = $values[$addr]
    $ws.Range($addr).Formula = "=1-" + $orig
}
